# Applies strikethrough formatting to the "old" project description
# (the Sessions / Example / User-registration walkthrough block), and
# splits a couple of runs so that Word's "session for all pages..." and
# the repeated "....." ellipsis get their own <w:proofErr> grammar
# markers, matching the target OOXML.

$d = $word.ActiveDocument

$ELLIPSIS = [char]0x2026

# 1) Simple paragraphs: just strike the whole paragraph (this also
#    adds the <w:pPr><w:rPr><w:strike/></w:rPr></w:pPr> to the
#    paragraph mark, and <w:rPr><w:strike/></w:rPr> to every run in
#    it -- exactly matching the diff).
foreach ($idx in 6, 8, 9, 10, 11, 12, 14, 15) {
    $para = $d.Paragraphs.Item($idx)
    $para.Range.Font.StrikeThrough = 1
}

# 2) "(session for all pages....)" paragraph: strike the paragraph,
#    then split the single run into "(" / "session" / " for all
#    pages....)" with proofErr gramStart/gramEnd bracketing "session".
$para = $d.Paragraphs.Item(13)
$para.Range.Font.StrikeThrough = 1
$start = $para.Range.Start
$end = $para.Range.End
$body = $d.Range($start, $end - 1)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' +
    '<w:r><w:rPr><w:strike/></w:rPr><w:t>(</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:strike/></w:rPr><w:t>session</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> for all pages' + $ELLIPSIS + '.)</w:t></w:r>' +
    '</w:p></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$body.InsertXML($xml)

# 3) "I should be able to see ... similar tastes....." paragraph:
#    strike the whole paragraph (adds strike to the existing runs and
#    the pPr), then split the trailing "....." off the last run into
#    its own proofErr-wrapped run.
$para = $d.Paragraphs.Item(18)
$para.Range.Font.StrikeThrough = 1
$start = $para.Range.Start
$end = $para.Range.End
$tail = $d.Range($end - 1 - 3, $end - 1)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:strike/></w:rPr><w:t>' + $ELLIPSIS + '..</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '</w:p></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$tail.InsertXML($xml)

Write-Host "applied strikethrough edits"
